# Update "想去人数" (want-to-go count) figures on a scrape re-run.
# 展览 (sheet1) and 全部类型 (sheet4) share the same 展览 rows; 演出 (sheet2)
# rows are also mirrored into 全部类型.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5736
$ws1.Range("F3").Value = 867
$ws1.Range("F4").Value = 85
$ws1.Range("F5").Value = 404
$ws1.Range("F6").Value = 8

# --- 演出 sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 62
$ws2.Range("F4").Value = 1

# --- 全部类型 sheet (combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5736
$ws4.Range("F3").Value = 867
$ws4.Range("F4").Value = 85
$ws4.Range("F5").Value = 62
$ws4.Range("F6").Value = 404
$ws4.Range("F7").Value = 8
$ws4.Range("F11").Value = 1
